# Update Pr(>F) values in column F (p-values rounded to 3 digits after
# removal of the y-axis log(x+1) scale from the biomass data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = 0.7554
    "F4"  = 0.7074
    "F7"  = 0.8939
    "F8"  = 0.731
    "F9"  = 0.4756
    "F12" = 0.7794
    "F13" = 0.7402
    "F14" = 0.1598
    "F17" = 0.9882
    "F18" = 0.2101
    "F19" = 0.7577
    "F22" = 0.0813
    "F23" = 0.0469
    "F24" = 0.1376
    "F27" = 0.8397
    "F28" = 0.7389
    "F29" = 0.6384
    "F32" = 0.8928
    "F33" = 0.0004
    "F34" = 0.9191
    "F37" = 0.6966
    "F38" = 0.0466
    "F39" = 0.7358
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
